$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Total Renda Fixa: Ago (I), Media (N) and Total (O) updated
$ws.Range("I2").Value = 300
$ws.Range("N2").Value = 552
$ws.Range("O2").Value = 6623.97

# Row 3 - CDB Banco Pan: Ago (I), Media (N) and Total (O) updated
$ws.Range("I3").Value = 300
$ws.Range("N3").Value = 311.82
$ws.Range("O3").Value = 3741.83

# Row 9 - BTG Fixa: Ago (I), Media (N) and Total (O) updated
$ws.Range("I9").Value = 300
$ws.Range("N9").Value = 137.55000000000001
$ws.Range("O9").Value = 1650.62

# Update the active selection to reflect the final cursor position
$ws.Activate()
$ws.Range("H16").Select()
